$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2: ECs -> Resolving-Mac (new TPM values) ---
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 15.428109
$ws.Range("H2").Value = 46.284327
$ws.Range("I2").Value = 0.105145687357564
$ws.Range("J2").Value = 0.105145687357564
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.888791333333333
$ws.Range("N2").Value = 5.666374
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 29.140478568922
$ws.Range("R2").Value = 262.264307120298
$ws.Range("S2").Value = 0.105145687357564
$ws.Range("T2").Value = 0.105145687357564

# --- Update row 3: FAPs -> Resolving-Mac (new TPM values) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 72.11798333333333
$ws.Range("H3").Value = 216.35395
$ws.Range("I3").Value = 0.4914986618531588
$ws.Range("J3").Value = 0.4914986618531588
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.888791333333333
$ws.Range("N3").Value = 5.666374
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 136.2158218974778
$ws.Range("R3").Value = 1225.9423970773
$ws.Range("S3").Value = 0.4914986618531588
$ws.Range("T3").Value = 0.4914986618531588

# --- Update row 4: MuSCs -> Resolving-Mac (new TPM values) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 17.753286
$ws.Range("H4").Value = 53.25985799999999
$ws.Range("I4").Value = 0.120992239510715
$ws.Range("J4").Value = 0.120992239510715
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.888791333333333
$ws.Range("N4").Value = 5.666374
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 33.532252734988
$ws.Range("R4").Value = 301.790274614892
$ws.Range("S4").Value = 0.120992239510715
$ws.Range("T4").Value = 0.120992239510715

# --- Update row 5: Resolving-Mac -> Resolving-Mac (new TPM values) ---
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 41.43140433333334
$ws.Range("H5").Value = 124.294213
$ws.Range("I5").Value = 0.2823634112785623
$ws.Range("J5").Value = 0.2823634112785622
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.888791333333333
$ws.Range("N5").Value = 5.666374
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 78.25527743262913
$ws.Range("R5").Value = 704.2974968936621
$ws.Range("S5").Value = 0.2823634112785623
$ws.Range("T5").Value = 0.2823634112785622

# --- Remove the now-obsolete rows 6-9 (old ECs target-cluster duplicate rows) ---
$ws.Range("A6:T9").EntireRow.Delete()
